# This workbook is a weekly "Poroto granado" (Macroferia Regional de Talca)
# price log. The weekly refresh inserts two new daily observations into the
# historical table (pushing every later row down), rather than just
# appending at the bottom.
#
#   - A new row is inserted at row 37 (date 44614 / 2022-02-22).
#   - A new row is inserted at row 113 of the (already-once-shifted) sheet
#     (date 44615 / 2022-02-23), which lands right after what was
#     previously row 112.
#
# Net effect: dimension grows from A1:R116 to A1:R118.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert #1: new row at 37 ---------------------------------------------
$ws.Rows(37).Insert()

$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 44614
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = 100112030
$ws.Range("G37").Value = "Poroto granado"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = 20000
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 800
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"

# --- Insert #2: new row at 113 (post first shift) --------------------------
$ws.Rows(113).Insert()

$ws.Range("A113").Value = 5
$ws.Range("B113").Value = "Macroferia Regional de Talca"
$ws.Range("C113").Value = "Maule"
$ws.Range("D113").Value = 44615
$ws.Range("E113").Value = 7
$ws.Range("F113").Value = 100112030
$ws.Range("G113").Value = "Poroto granado"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 20000
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = 20000
$ws.Range("N113").Value = "$/saco 25 kilos"
$ws.Range("O113").Value = "Región del Maule"
$ws.Range("P113").Value = 800
$ws.Range("Q113").Value = 25
$ws.Range("R113").Value = "Hortaliza"
